# Fruta / hortaliza, semanal
# Inserts a new weekly price record for "Ajo" at row 223 of Sheet1,
# shifting the existing rows 223:235 down to 224:236.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 223; this pushes the current
# rows 223-235 down to 224-236 and extends the used range to R236.
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new record's data.
$ws.Cells.Item(223, 1).Value2 = 5
$ws.Cells.Item(223, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(223, 3).Value2 = "Maule"
$ws.Cells.Item(223, 4).Value2 = 44585
$ws.Cells.Item(223, 5).Value2 = 7
$ws.Cells.Item(223, 6).Value2 = 100112003
$ws.Cells.Item(223, 7).Value2 = "Ajo"
$ws.Cells.Item(223, 8).Value2 = "Chino"
$ws.Cells.Item(223, 9).Value2 = "Primera"
$ws.Cells.Item(223, 10).Value2 = 200
$ws.Cells.Item(223, 11).Value2 = 18000
$ws.Cells.Item(223, 12).Value2 = 18000
$ws.Cells.Item(223, 13).Value2 = 18000
$ws.Cells.Item(223, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(223, 15).Value2 = "China"
$ws.Cells.Item(223, 16).Value2 = 1800
$ws.Cells.Item(223, 17).Value2 = 10
$ws.Cells.Item(223, 18).Value2 = "Hortaliza"
